$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# This edit inserts a new "2022-Q1" sheet (same 7-column "fund" layout
# as "2021-Q4") right before the "总计" summary sheet, and prepends a
# matching row to "总计".
#
# The target sheetId numbering (2022-Q1 -> 3, 总计 -> 4, i.e. 总计 gets
# a brand-new id while the freed-up old id 3 is reused by the new
# sheet) only falls out naturally if "总计" is deleted and rebuilt, so
# that is the approach used below; its data/format are easy to
# reconstruct since the sheet only has a handful of cells.
# ---------------------------------------------------------------------

$totalOld = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

[void]$totalOld.Delete()

# New "2022-Q1" sheet: duplicate "2021-Q4" (identical headers/style) and
# land right after it (end of workbook at this point).
$q4.Copy([System.Type]::Missing, [System.Type]::Missing)
$q1Sheet = $wb.Worksheets.Item("2021-Q4 (2)")
$q1Sheet.Name = "2022-Q1"

# Row 2 already has 539002 / 建信新兴市场优选混合QDII in B2/C2 (unchanged).
# D2:G2 are stored as text in the source file, so force text entry via a
# leading apostrophe and then strip the resulting quote-prefix style back
# off (restore the plain, un-styled look the other data cells use).
$plainStyle = $q1Sheet.Range("B2").Style

$q1Sheet.Range("D2").Value = "'0.14"
$q1Sheet.Range("E2").Value = "'83.76"
$q1Sheet.Range("F2").Value = "'3.77"
$q1Sheet.Range("G2").Value = "'0.0053"
$q1Sheet.Range("D2:G2").Style = $plainStyle

$q1Sheet.Range("H2").Value = 10

# New "总计" sheet: duplicate "2021-Q4" again (to land right after
# "2022-Q1" and to pick up the next sheetId), then rebuild it as the
# 4-column summary sheet.
$q4.Copy([System.Type]::Missing, [System.Type]::Missing)
$totalNew = $wb.Worksheets.Item("2021-Q4 (2)")
$totalNew.Name = "总计"

$totalNew.Cells.Clear()

# Reuse the header / index-column formatting (border+bold+center for
# row 1, center-top for column A) from the sheet we just built.
$q1Sheet.Range("B1").Copy()
$totalNew.Range("B1:D1").PasteSpecial(-4122)

$q1Sheet.Range("A2").Copy()
$totalNew.Range("A2:A4").PasteSpecial(-4122)

$totalNew.Range("B1").Value = "日期"
$totalNew.Range("C1").Value = "持有数量(只)"
$totalNew.Range("D1").Value = "持有市值(亿元)"

$totalNew.Range("A2").Value = 0
$totalNew.Range("B2").Value = "2022-Q1"
$totalNew.Range("C2").Value = 1
$totalNew.Range("D2").Value = 0.01

$totalNew.Range("A3").Value = 1
$totalNew.Range("B3").Value = "2021-Q4"
$totalNew.Range("C3").Value = 1
$totalNew.Range("D3").Value = 0.01

$totalNew.Range("A4").Value = 2
$totalNew.Range("B4").Value = "2021-Q3"
$totalNew.Range("C4").Value = 1
$totalNew.Range("D4").Value = 0.01

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("2021-Q3").Activate()
